$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column F (dSF) values for rows 3-9
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 3
$ws.Range("F5").Value = -2
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -5
$ws.Range("F9").Value = -1
